$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.092.12"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "'1.829.47"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'243.05"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.07452"
$ws.Range("E8").Value = "  -1.77%  "
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("D10").Value = "'23.29"
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("D11").Value = "'0.07700"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").Value = "'1.832.19"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "'5.006"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").Value = "'0.6676"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "'82.47"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").Value = "'0.000009379"
$ws.Range("E16").Value = "  -6.48%  "
$ws.Range("D17").Value = "'5.956"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").Value = "'29.104.14"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "'2.082.57"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").Value = "'12.59"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").Value = "'223.09"
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'7.123"
$ws.Range("E23").Value = "  -1.50%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "'159.95"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "'0.1392"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "'8.493"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").Value = "'17.89"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").Value = "'1.492"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").Value = "'0.05785"
$ws.Range("E30").Value = "  +9.65%  "
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("D32").Value = "'4.123"
$ws.Range("E32").Value = "  +2.43%  "
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'1.828"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7383"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "'2.670"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "'1.225.66"
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.762"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").Value = "'6.490"
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("D42").Value = "'0.8921"
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'102.05"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "'0.00000000126"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "'1.972.09"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").Value = "'65.93"
$ws.Range("E47").Value = "  +2.36%  "
$ws.Range("D48").Value = "'0.5087"
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("D49").Value = "'0.07583"
$ws.Range("E49").Value = "  +14.12%  "
$ws.Range("D50").Value = "'0.4058"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").Value = "'8.975"
$ws.Range("E51").Value = "  +0.40%  "
